$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.30506
$ws.Range("H2").Value = 0.91518
$ws.Range("I2").Value = 0.07720017721337037
$ws.Range("J2").Value = 0.07720017721337037
$ws.Range("M2").Value = 3.073545
$ws.Range("N2").Value = 9.220635
$ws.Range("O2").Value = 0.6324804786546022
$ws.Range("P2").Value = 0.6324804786546022
$ws.Range("Q2").Value = 0.9376156376999999
$ws.Range("R2").Value = 8.4385407393
$ws.Range("S2").Value = 0.04882760503613261
$ws.Range("T2").Value = 0.04882760503613261

# Row 3
$ws.Range("G3").Value = 0.30506
$ws.Range("H3").Value = 0.91518
$ws.Range("I3").Value = 0.07720017721337037
$ws.Range("J3").Value = 0.07720017721337037
$ws.Range("O3").Value = 0.215691788931517
$ws.Range("P3").Value = 0.215691788931517
$ws.Range("Q3").Value = 0.3197505710466667
$ws.Range("R3").Value = 2.87775513942
$ws.Range("S3").Value = 0.01665144432898199
$ws.Range("T3").Value = 0.01665144432898199

# Row 4
$ws.Range("G4").Value = 0.30506
$ws.Range("H4").Value = 0.91518
$ws.Range("I4").Value = 0.07720017721337037
$ws.Range("J4").Value = 0.07720017721337037
$ws.Range("O4").Value = 0.1518277324138807
$ws.Range("P4").Value = 0.1518277324138807
$ws.Range("Q4").Value = 0.2250758101666666
$ws.Range("R4").Value = 2.0256822915
$ws.Range("S4").Value = 0.01172112784825577
$ws.Range("T4").Value = 0.01172112784825577

# Row 5
$ws.Range("G5").Value = 2.879048666666666
$ws.Range("H5").Value = 8.637146
$ws.Range("I5").Value = 0.7285880393122151
$ws.Range("J5").Value = 0.7285880393122151
$ws.Range("M5").Value = 3.073545
$ws.Range("N5").Value = 9.220635
$ws.Range("O5").Value = 0.6324804786546022
$ws.Range("P5").Value = 0.6324804786546022
$ws.Range("Q5").Value = 8.848885634189998
$ws.Range("R5").Value = 79.63997070770999
$ws.Range("S5").Value = 0.4608177118462079
$ws.Range("T5").Value = 0.4608177118462079

# Row 6
$ws.Range("G6").Value = 2.879048666666666
$ws.Range("H6").Value = 8.637146
$ws.Range("I6").Value = 0.7285880393122151
$ws.Range("J6").Value = 0.7285880393122151
$ws.Range("O6").Value = 0.215691788931517
$ws.Range("P6").Value = 0.215691788931517
$ws.Range("Q6").Value = 3.017693093941555
$ws.Range("S6").Value = 0.1571504575933581
$ws.Range("T6").Value = 0.1571504575933581

# Row 7
$ws.Range("G7").Value = 2.879048666666666
$ws.Range("H7").Value = 8.637146
$ws.Range("I7").Value = 0.7285880393122151
$ws.Range("J7").Value = 0.7285880393122151
$ws.Range("O7").Value = 0.1518277324138807
$ws.Range("P7").Value = 0.1518277324138807
$ws.Range("Q7").Value = 2.124186098338888
$ws.Range("S7").Value = 0.110619869872649
$ws.Range("T7").Value = 0.110619869872649

# Row 8
$ws.Range("G8").Value = 0.7674366666666668
$ws.Range("I8").Value = 0.1942117834744146
$ws.Range("J8").Value = 0.1942117834744146
$ws.Range("M8").Value = 3.073545
$ws.Range("N8").Value = 9.220635
$ws.Range("O8").Value = 0.6324804786546022
$ws.Range("P8").Value = 0.6324804786546022
$ws.Range("Q8").Value = 2.35875112965
$ws.Range("R8").Value = 21.22876016685
$ws.Range("S8").Value = 0.1228351617722617
$ws.Range("T8").Value = 0.1228351617722617

# Row 9
$ws.Range("G9").Value = 0.7674366666666668
$ws.Range("I9").Value = 0.1942117834744146
$ws.Range("J9").Value = 0.1942117834744146
$ws.Range("O9").Value = 0.215691788931517
$ws.Range("P9").Value = 0.215691788931517
$ws.Range("Q9").Value = 0.8043936025988889
$ws.Range("R9").Value = 7.239542423390001
$ws.Range("S9").Value = 0.04188988700917691
$ws.Range("T9").Value = 0.04188988700917692

# Row 10
$ws.Range("G10").Value = 0.7674366666666668
$ws.Range("I10").Value = 0.1942117834744146
$ws.Range("J10").Value = 0.1942117834744146
$ws.Range("O10").Value = 0.1518277324138807
$ws.Range("P10").Value = 0.1518277324138807
$ws.Range("R10").Value = 5.09599051175
$ws.Range("S10").Value = 0.02948673469297595
$ws.Range("T10").Value = 0.02948673469297596
